# Add a new final slide ("СПАСИБО ЗА ВНИМАНИЕ" / "Thank you for your attention!")
# using the "Title Only" layout (same layout family as PowerPoint's ppLayoutTitleOnly),
# matching ppt/slideLayouts/slideLayout6.xml ("Только заголовок").

$p = $ppt.ActivePresentation

# Insert as the new (13th) slide, right after the current last slide.
$s = $p.Slides.Add($p.Slides.Count + 1, 11)

$shp = $s.Shapes.Item(1)
$shp.Name = "Заголовок 1"

# Position/size of the title placeholder (EMU -> points, with a tiny
# sub-EMU nudge so the points->EMU round trip lands back on the exact
# EMU value instead of being truncated one unit short).
$emuPerPt = 12700
$halfEmuPt = 0.5 / $emuPerPt

$shp.Left = 485745 / $emuPerPt + $halfEmuPt
$shp.Top = 2481943 / $emuPerPt + $halfEmuPt
$shp.Width = 10173546 / $emuPerPt + $halfEmuPt
$shp.Height = 1320800 / $emuPerPt + $halfEmuPt

# No auto-fit for the title text box.
$shp.TextFrame.AutoSize = 0

$tr = $shp.TextFrame.TextRange
$tr.Text = "Спасибо за внимание!"
$tr.Font.Size = 72
$tr.LanguageID = "ru-RU"
